$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new customer row (fetched from Excel data) below the existing
# username/password rows.
$ws.Range("A6").Value = "mngr434372"
$ws.Range("B6").Value = "umAsapE"

# Move the active selection to the newly added cell, matching the
# workbook state after the edit.
$ws.Range("B6").Select()
